$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "10/11" result for the Taekwondo position column (G) on the
# Dynamic time warping row, and reuse the existing "11/12" value for the
# Elastic action comparison with freedom degree row.
$ws.Range("G2").Value = "10/11"
$ws.Range("G3").Value = "11/12"

# Match the author's last active selection/cell.
$ws.Range("G2").Select()
